$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

# Insert a new column before column A; existing columns (A-E) shift to (B-F).
$ws.Columns("A").Insert()

# Populate the new "Identificador" id column.
$ws.Range("A1").Value = "Identificador"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2

# Match the header's bold style used by the rest of row 1.
$ws.Range("A1").Font.Bold = $true

# Column A ends up auto-fit to its content at 12 characters wide, same as the
# other header columns in this sheet. Excel's VBA ColumnWidth (character
# units) maps to an on-disk width that is 5/6 larger than the assigned value,
# so back that offset out to land exactly on a stored width of 12.
$ws.Columns("A").ColumnWidth = 11.166666666666666
